$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 4.918060118371062
$ws.Cells.Item(2, 4).Value = 4.846278250066219
$ws.Cells.Item(2, 5).Value = 10.4224791904717
$ws.Cells.Item(2, 6).Value = 27.58103282242811
$ws.Cells.Item(2, 7).Value = 35.16215282412897
$ws.Cells.Item(2, 8).Value = 15.31823998826257
$ws.Cells.Item(2, 11).Value = 13.90623777669959
$ws.Cells.Item(2, 13).Value = 15.66677161891834
$ws.Cells.Item(2, 14).Value = 17.55335470820703

$ws.Cells.Item(3, 3).Value = 4.749848789405855
$ws.Cells.Item(3, 4).Value = 4.861085847570143
$ws.Cells.Item(3, 5).Value = 10.3262835355799
$ws.Cells.Item(3, 6).Value = 27.40769256138062
$ws.Cells.Item(3, 7).Value = 34.77530741027384
$ws.Cells.Item(3, 8).Value = 15.3301704135749
$ws.Cells.Item(3, 11).Value = 13.41446731855034
$ws.Cells.Item(3, 13).Value = 15.38426238701327
$ws.Cells.Item(3, 14).Value = 17.62677967540915

$ws.Cells.Item(4, 3).Value = 4.645085964189617
$ws.Cells.Item(4, 4).Value = 4.870542902905159
$ws.Cells.Item(4, 5).Value = 10.27055442675591
$ws.Cells.Item(4, 6).Value = 27.31074594068393
$ws.Cells.Item(4, 7).Value = 34.55080680784014
$ws.Cells.Item(4, 8).Value = 15.34151291110592
$ws.Cells.Item(4, 11).Value = 13.10631007539394
$ws.Cells.Item(4, 13).Value = 15.21250582492889
$ws.Cells.Item(4, 14).Value = 17.67379336881909

$ws.Cells.Item(5, 3).Value = 4.602098494846431
$ws.Cells.Item(5, 4).Value = 4.874488994477331
$ws.Cells.Item(5, 5).Value = 10.24870439190044
$ws.Cells.Item(5, 6).Value = 27.27365587642398
$ws.Cells.Item(5, 7).Value = 34.46270145193209
$ws.Cells.Item(5, 8).Value = 15.3471419013325
$ws.Cells.Item(5, 11).Value = 12.97938073577536
$ws.Cells.Item(5, 13).Value = 15.1430402625757
$ws.Cells.Item(5, 14).Value = 17.69343915041156

$ws.Cells.Item(6, 3).Value = 4.594944943723196
$ws.Cells.Item(6, 4).Value = 4.875149825788967
$ws.Cells.Item(6, 5).Value = 10.24512873593673
$ws.Cells.Item(6, 6).Value = 27.26764387017797
$ws.Cells.Item(6, 7).Value = 34.44827869699716
$ws.Cells.Item(6, 8).Value = 15.34813729417574
$ws.Cells.Item(6, 11).Value = 12.95822896125779
$ws.Cells.Item(6, 13).Value = 15.131540225698
$ws.Cells.Item(6, 14).Value = 17.69673080580136

$ws.Cells.Item(7, 3).Value = 4.644507312309413
$ws.Cells.Item(7, 4).Value = 4.870595747129096
$ws.Cells.Item(7, 5).Value = 10.27025624070575
$ws.Cells.Item(7, 6).Value = 27.3102359091823
$ws.Cells.Item(7, 7).Value = 34.54960476847898
$ws.Cells.Item(7, 8).Value = 15.34158475385554
$ws.Cells.Item(7, 11).Value = 13.10460345040279
$ws.Cells.Item(7, 13).Value = 15.21156672647668
$ws.Cells.Item(7, 14).Value = 17.67405634310597

$ws.Cells.Item(8, 3).Value = 4.8604121719033
$ws.Cells.Item(8, 4).Value = 4.851308397552581
$ws.Cells.Item(8, 5).Value = 10.3886320738168
$ws.Cells.Item(8, 6).Value = 27.51931679162903
$ws.Cells.Item(8, 7).Value = 35.02613141707781
$ws.Cells.Item(8, 8).Value = 15.32151798483817
$ws.Cells.Item(8, 11).Value = 13.73808311523981
$ws.Cells.Item(8, 13).Value = 15.5690671195614
$ws.Cells.Item(8, 14).Value = 17.57827211119905

$ws.Cells.Item(9, 3).Value = 5.268998447748573
$ws.Cells.Item(9, 4).Value = 4.816362523010491
$ws.Cells.Item(9, 5).Value = 10.6461989714604
$ws.Cells.Item(9, 6).Value = 28.00297431940131
$ws.Cells.Item(9, 7).Value = 36.05887763048881
$ws.Cells.Item(9, 8).Value = 15.31416911072383
$ws.Cells.Item(9, 11).Value = 14.92286128474466
$ws.Cells.Item(9, 13).Value = 16.27938745495121
$ws.Cells.Item(9, 14).Value = 17.4056692413071

$ws.Cells.Item(10, 3).Value = 5.556504302873424
$ws.Cells.Item(10, 4).Value = 4.792412498579704
$ws.Cells.Item(10, 5).Value = 10.84943419672811
$ws.Cells.Item(10, 6).Value = 28.40075469833503
$ws.Cells.Item(10, 7).Value = 36.87026443463502
$ws.Cells.Item(10, 8).Value = 15.32842262359443
$ws.Cells.Item(10, 11).Value = 15.74857206094275
$ws.Cells.Item(10, 13).Value = 16.80130309770293
$ws.Cells.Item(10, 14).Value = 17.28801746661947

$ws.Cells.Item(11, 3).Value = 5.683893119085585
$ws.Cells.Item(11, 4).Value = 4.781885392800195
$ws.Cells.Item(11, 5).Value = 10.94458263588374
$ws.Cells.Item(11, 6).Value = 28.5903227791195
$ws.Cells.Item(11, 7).Value = 37.24915439943913
$ws.Cells.Item(11, 8).Value = 15.33919362078088
$ws.Cells.Item(11, 11).Value = 16.11279353088346
$ws.Cells.Item(11, 13).Value = 17.03760952377139
$ws.Cells.Item(11, 14).Value = 17.23645671617937

$ws.Cells.Item(12, 3).Value = 5.731596766687468
$ws.Cells.Item(12, 4).Value = 4.777951494647026
$ws.Cells.Item(12, 5).Value = 10.98096981676967
$ws.Cells.Item(12, 6).Value = 28.66329017045752
$ws.Cells.Item(12, 7).Value = 37.39389408993371
$ws.Cells.Item(12, 8).Value = 15.34388934263996
$ws.Cells.Item(12, 11).Value = 16.24895588206696
$ws.Cells.Item(12, 13).Value = 17.12684226421648
$ws.Cells.Item(12, 14).Value = 17.21721177855758

$ws.Cells.Item(13, 3).Value = 5.721347449564375
$ws.Cells.Item(13, 4).Value = 4.778796402065455
$ws.Cells.Item(13, 5).Value = 10.97311783142761
$ws.Cells.Item(13, 6).Value = 28.64752369906279
$ws.Cells.Item(13, 7).Value = 37.36266798835999
$ws.Cells.Item(13, 8).Value = 15.34285059281096
$ws.Cells.Item(13, 11).Value = 16.21971095229159
$ws.Cells.Item(13, 13).Value = 17.10763697975711
$ws.Cells.Item(13, 14).Value = 17.22134409562102

$ws.Cells.Item(14, 3).Value = 5.687828713414621
$ws.Cells.Item(14, 4).Value = 4.781560699030995
$ws.Cells.Item(14, 5).Value = 10.94756926615807
$ws.Cells.Item(14, 6).Value = 28.5963024954516
$ws.Cells.Item(14, 7).Value = 37.2610377575967
$ws.Cells.Item(14, 8).Value = 15.33956757466522
$ws.Cells.Item(14, 11).Value = 16.12403161279597
$ws.Cells.Item(14, 13).Value = 17.04495629879867
$ws.Cells.Item(14, 14).Value = 17.23486782222013

$ws.Cells.Item(15, 3).Value = 5.66722645086167
$ws.Cells.Item(15, 4).Value = 4.783260734124763
$ws.Cells.Item(15, 5).Value = 10.9319655289586
$ws.Cells.Item(15, 6).Value = 28.56508023441616
$ws.Cells.Item(15, 7).Value = 37.19894636889136
$ws.Cells.Item(15, 8).Value = 15.33763698510955
$ws.Cells.Item(15, 11).Value = 16.06519265095537
$ws.Cells.Item(15, 13).Value = 17.00652720843675
$ws.Cells.Item(15, 14).Value = 17.24318790635256

$ws.Cells.Item(16, 3).Value = 5.548106607702352
$ws.Cells.Item(16, 4).Value = 4.793107837005041
$ws.Cells.Item(16, 5).Value = 10.84326762346172
$ws.Cells.Item(16, 6).Value = 28.38853455465982
$ws.Cells.Item(16, 7).Value = 36.84568764977335
$ws.Cells.Item(16, 8).Value = 15.32780503660494
$ws.Cells.Item(16, 11).Value = 15.72452924870034
$ws.Cells.Item(16, 13).Value = 16.78582959786464
$ws.Cells.Item(16, 14).Value = 17.29142636227819

$ws.Cells.Item(17, 3).Value = 5.474124704882811
$ws.Cells.Item(17, 4).Value = 4.799242643361967
$ws.Cells.Item(17, 5).Value = 10.78952286119713
$ws.Cells.Item(17, 6).Value = 28.28239667807765
$ws.Cells.Item(17, 7).Value = 36.63137871840366
$ws.Cells.Item(17, 8).Value = 15.32287201399995
$ws.Cells.Item(17, 11).Value = 15.51253270365688
$ws.Cells.Item(17, 13).Value = 16.65008700323837
$ws.Cells.Item(17, 14).Value = 17.3215197231174

$ws.Cells.Item(18, 3).Value = 5.431254783392384
$ws.Cells.Item(18, 4).Value = 4.802805872673582
$ws.Cells.Item(18, 5).Value = 10.75886594347618
$ws.Cells.Item(18, 6).Value = 28.22216354039489
$ws.Cells.Item(18, 7).Value = 36.50904381641062
$ws.Cells.Item(18, 8).Value = 15.32043819017195
$ws.Cells.Item(18, 11).Value = 15.38953126832465
$ws.Cells.Item(18, 13).Value = 16.57191163578646
$ws.Cells.Item(18, 14).Value = 17.33901318497401

$ws.Cells.Item(19, 3).Value = 5.416686756687486
$ws.Cells.Item(19, 4).Value = 4.804018284876547
$ws.Cells.Item(19, 5).Value = 10.74853088210772
$ws.Cells.Item(19, 6).Value = 28.2019112217074
$ws.Cells.Item(19, 7).Value = 36.46778742614489
$ws.Cells.Item(19, 8).Value = 15.31968341873235
$ws.Cells.Item(19, 11).Value = 15.34770584606595
$ws.Cells.Item(19, 13).Value = 16.54542849518621
$ws.Cells.Item(19, 14).Value = 17.34496792307615

$ws.Cells.Item(20, 3).Value = 5.482033420955638
$ws.Cells.Item(20, 4).Value = 4.798585999564301
$ws.Cells.Item(20, 5).Value = 10.79521784751737
$ws.Cells.Item(20, 6).Value = 28.29361133531692
$ws.Cells.Item(20, 7).Value = 36.65409707160753
$ws.Cells.Item(20, 8).Value = 15.32335537440123
$ws.Cells.Item(20, 11).Value = 15.53521141649185
$ws.Cells.Item(20, 13).Value = 16.66454803707562
$ws.Cells.Item(20, 14).Value = 17.3182971456886

$ws.Cells.Item(21, 3).Value = 5.69768886564934
$ws.Cells.Item(21, 4).Value = 4.780747336742659
$ws.Cells.Item(21, 5).Value = 10.95506407396816
$ws.Cells.Item(21, 6).Value = 28.61131579045714
$ws.Cells.Item(21, 7).Value = 37.29085592334903
$ws.Cells.Item(21, 8).Value = 15.34051513092567
$ws.Cells.Item(21, 11).Value = 16.1521836243185
$ws.Cells.Item(21, 13).Value = 17.06337466018917
$ws.Cells.Item(21, 14).Value = 17.23088799104361

$ws.Cells.Item(22, 3).Value = 5.835491480263432
$ws.Cells.Item(22, 4).Value = 4.769394494307705
$ws.Cells.Item(22, 5).Value = 11.06159614359335
$ws.Cells.Item(22, 6).Value = 28.82581580045433
$ws.Cells.Item(22, 7).Value = 37.71430433737811
$ws.Cells.Item(22, 8).Value = 15.35532606282566
$ws.Cells.Item(22, 11).Value = 16.54510356566198
$ws.Cells.Item(22, 13).Value = 17.32252659591903
$ws.Cells.Item(22, 14).Value = 17.17539228001872

$ws.Cells.Item(23, 3).Value = 5.76224513090451
$ws.Cells.Item(23, 4).Value = 4.775425878123632
$ws.Cells.Item(23, 5).Value = 11.00455960756241
$ws.Cells.Item(23, 6).Value = 28.71072473928432
$ws.Cells.Item(23, 7).Value = 37.48768309538645
$ws.Cells.Item(23, 8).Value = 15.34709213761772
$ws.Cells.Item(23, 11).Value = 16.33637418294985
$ws.Cells.Item(23, 13).Value = 17.18437890836555
$ws.Cells.Item(23, 14).Value = 17.20486271888615

$ws.Cells.Item(24, 3).Value = 5.478458935698528
$ws.Cells.Item(24, 4).Value = 4.79888275526585
$ws.Cells.Item(24, 5).Value = 10.79264238912371
$ws.Cells.Item(24, 6).Value = 28.28853873322096
$ws.Cells.Item(24, 7).Value = 36.64382337019521
$ws.Cells.Item(24, 8).Value = 15.32313559424926
$ws.Cells.Item(24, 11).Value = 15.5249618551976
$ws.Cells.Item(24, 13).Value = 16.65801061764916
$ws.Cells.Item(24, 14).Value = 17.31975347350926

$ws.Cells.Item(25, 3).Value = 5.160451390741672
$ws.Cells.Item(25, 4).Value = 4.825511363026762
$ws.Cells.Item(25, 5).Value = 10.57394689530474
$ws.Cells.Item(25, 6).Value = 27.86448515103269
$ws.Cells.Item(25, 7).Value = 35.76970001013813
$ws.Cells.Item(25, 8).Value = 15.31271414234511
$ws.Cells.Item(25, 11).Value = 14.6095723285995
$ws.Cells.Item(25, 13).Value = 16.08682285160658
$ws.Cells.Item(25, 14).Value = 17.45074538275109
